$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column as text so numeric-looking values (e.g. "1.003")
# are not coerced into numbers by Excel, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.171.01"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "1.800.50"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "324.21"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "0.4294"
$ws.Range("E7").Value = "  -3.21%  "
$ws.Range("D8").Value = "0.3629"
$ws.Range("E8").Value = "  -3.00%  "
$ws.Range("D9").Value = "44.74"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "0.07557"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "1.123"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "21.71"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "6.201"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").Value = "7.371"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "1.814.80"
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("D17").Value = "92.86"
$ws.Range("E17").Value = "  +4.48%  "
$ws.Range("D18").Value = "0.00001069"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "0.06350"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").Value = "1.002"
$ws.Range("D21").Value = "17.26"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "5.990"
$ws.Range("E22").Value = "  -3.14%  "
$ws.Range("D23").Value = "28.175.35"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").Value = "11.42"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").Value = "2.164"
$ws.Range("E25").Value = "  -6.61%  "
$ws.Range("D26").Value = "159.67"
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("D27").Value = "20.43"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "2.021.23"
$ws.Range("E28").Value = "  +3.32%  "
$ws.Range("D29").Value = "2.239"
$ws.Range("E29").Value = "  -5.11%  "
$ws.Range("D30").Value = "128.18"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "1.178"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").Value = "5.896"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").Value = "0.09030"
$ws.Range("E33").Value = "  -3.49%  "
$ws.Range("D34").Value = "3.527"
$ws.Range("E34").Value = "  -3.39%  "
$ws.Range("D35").Value = "12.86"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("D36").Value = "0.02362"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "5.143"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").Value = "0.6511"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "0.06130"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "0.2126"
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("D41").Value = "1.194"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "1.432"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "7.975"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "13.55"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").Value = "0.6020"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "3.711"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "125.48"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").Value = "1.988"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "1.159"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "0.06970"
$ws.Range("E51").Value = "  +1.06%  "
